$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$newRows = @(
    @('22CRB00136', 'Bunner', 'DOMESTIC VIOLENCE', '2919.25(A)', 'No Data', 'No Contest', 'Guilty', '$ 50', '$ 25', 'None', 'None'),
    @('22CRB00136', 'Bunner', 'ASSAULT - M1', '2903.13(A)', 'No Data', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None'),
    @('22CRB00136', 'Bunner', 'DOMESTIC VIOLENCE', '2919.25(A)', 'No Data', 'No Contest', 'Guilty', '$ 50', '$ 25', 'None', 'None'),
    @('22CRB00136', 'Bunner', 'ASSAULT - M1', '2903.13(A)', 'No Data', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None'),
    @('22CRB00136', 'Hemmeter', 'DOMESTIC VIOLENCE', '2919.25(A)', 'No Data', 'No Contest', 'Guilty', '$ 0', '$ 0', '5', 'None'),
    @('22CRB00136', 'Hemmeter', 'ASSAULT - M1', '2903.13(A)', 'No Data', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None'),
    @('22CRB00136', 'Hemmeter', 'DOMESTIC VIOLENCE', '2919.25(A)', 'No Data', 'No Contest', 'Guilty', '$ 50', '$ 25', '2', 'None'),
    @('22CRB00136', 'Hemmeter', 'ASSAULT - M1', '2903.13(A)', 'No Data', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None'),
    @('22CRB00136', 'Hemmeter', 'DOMESTIC VIOLENCE', '2919.25(A)', 'No Data', 'No Contest', 'Guilty', '$ 50', '$ 25', '2', 'None'),
    @('22CRB00136', 'Hemmeter', 'ASSAULT - M1', '2903.13(A)', 'No Data', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None'),
    @('22CRB00136', 'Hemmeter', 'DOMESTIC VIOLENCE', '2919.25(A)', 'No Data', 'No Contest', 'Guilty', '$ 50', '$ 25', '10', 'None'),
    @('22CRB00136', 'Hemmeter', 'ASSAULT - M1', '2903.13(A)', 'No Data', 'No Contest', 'Guilty', '$ 0', '$ 0', 'None', 'None'),
)

$startRow = 778
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        Set-TextValue $ws.Cells.Item($r, $c + 1) $rowData[$c]
    }
}

